$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the expiration date for the UT row (row 2)
$ws.Range("E2").Value = 99999999

# Remove the now-obsolete row (old row 3, the SS-test VIN row)
$ws.Rows("3:3").Delete()

# Select the row below the data, matching the post-edit selection state
$ws.Range("A3:XFD8").Select()
